$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at the very top; all existing data (previously rows 1-12)
# shifts down to rows 2-13.
$ws.Rows.Item(1).Insert()

# Populate the new header row: A1 is a blank/empty text label, B1:D1 hold
# the column headers "X", "Y", "Z".
# A single apostrophe forces a literal-text entry of empty content (Excel's
# "quote prefix" convention), which lands in the shared-string table as an
# empty string rather than clearing the cell outright.
$ws.Range("A1").Value = "'"
$ws.Range("B1").Value = "X"
$ws.Range("C1").Value = "Y"
$ws.Range("D1").Value = "Z"

# The quote-prefix entry above also stamps A1 with a "quotePrefix" cell
# style; strip that back out by pasting in the (default) formatting from a
# neighboring cell, so A1 keeps its empty text but reverts to the normal
# style.
$ws.Range("B2").Copy()
$ws.Range("A1").PasteSpecial(-4122)
$excel.CutCopyMode = $false
